$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "along" in E1
$ws.Range("E1").Value = "along"

# New data values: A (label stays same), B (mean), C (CI_low), D (CI_high), E (new constant string)
$labels = @("<b>All</b>", "<b>Europe</b>", "France", "Germany", "Italy", "Poland", "Spain", "United Kingdom", "Switzerland", "Japan", "USA")
$mean   = @(-0.00909143460136564, -0.0248745755214356, 0.00380154502489296, -0.0324505655879783, -0.124474901990024, -0.0485006320675743, -0.144349327572778, 0.113637141813141, 0.114096038164518, 0.0351871074271903, -0.0240735508390279)
$cilow  = @(-0.0391044692896361, -0.0626813247689267, -0.112951374195182, -0.107147591128177, -0.213438715323915, -0.22655316621991, -0.234381326674561, 0.0240508085627794, -0.0103402651593815, -0.0215855864789808, -0.119671251686483)
$cihigh = @(0.0209216000869048, 0.0129321737260554, 0.120554464244968, 0.0422464599522202, -0.0355110886561329, 0.129551902084762, -0.0543173284709959, 0.203223475063502, 0.238532341488417, 0.0919598013333613, 0.0715241500084268)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]
    $ws.Cells.Item($row, 2).Value = $mean[$i]
    $ws.Cells.Item($row, 3).Value = $cilow[$i]
    $ws.Cells.Item($row, 4).Value = $cihigh[$i]
    $ws.Cells.Item($row, 5).Value = "cut_aid_in_programTRUE"
}
